$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: new "Social Logins" / Clerk entry -----------------------------
$ws.Range("A9").Value = "Social Logins"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = "Clerk"
$ws.Range("D9").Value = "`$25 after 10,000 monthly users"
$ws.Range("E9").Value = "10,000/ Users a month"
$ws.Range("F9").Value = "https://clerk.com/pricing"

# Hyperlink for the notes column, mirroring the other rows
$ws.Hyperlinks.Add($ws.Range("E9"), "https://clerk.com/pricing")

# Match formatting of the row above (this also restores the correct
# hyperlink-cell style on E9 after Hyperlinks.Add stomped on it)
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 13: Total row -----------------------------------------------------
$ws.Range("A13").Value = "Total"
$ws.Range("B13").Formula = "=SUM(B1:B9)"

$ws.Range("A13").Interior.ThemeColor = 3
$ws.Range("B13").Interior.Color = 49407

$ws.Range("C18").Select()
